$d = $word.ActiveDocument

function New-FlatOpcXml($bodyInnerXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Paragraph 1 (bold heading): split the single run into four runs so that
# "Successful" / " or " / "failed" / ")." are each their own run, dropping
# ", canceled or live" (keeping the wording "Successful or failed).").
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$xml1 = New-FlatOpcXml('<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Question 3: Compare the length of campaigns to outcomes (Successful</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> or </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>failed</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>).</w:t></w:r></w:p>')
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Paragraph 3 (plain single-run paragraph): simple text replacement.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This question has proven difficult to visualize, given time constraints, beyond the scope of this project. However, it would provide valuable insights. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Visualization of the data does not show a solid correlation between the length of the campaign and its success or failure.",
    2
)

# ---------------------------------------------------------------------------
# Paragraph 5: replace the single run with five runs.
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range
$xml5 = New-FlatOpcXml('<w:p><w:r><w:t xml:space="preserve">When all countries are viewed together it appears that success or failure </w:t></w:r><w:r><w:t>is spread</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>evenly</w:t></w:r><w:r><w:t xml:space="preserve"> across the number of campaign days.</w:t></w:r></w:p>')
$r5.InsertXML($xml5)

# ---------------------------------------------------------------------------
# Paragraph 7 ("This analysis should be completed for future work.") is
# replaced by two runs of new text, followed by a new empty paragraph and a
# new paragraph of text. Since paragraph 7 currently is the LAST paragraph in
# the document body, InsertXML with a full <w:p> wrapper there would append
# an extra paragraph instead of replacing in place (the trailing paragraph
# mark is tied to the section properties). To avoid that, first grow the
# body with plain placeholder paragraphs (so paragraph 7 is no longer last),
# then fill each paragraph's content in place, and finally trim the spare
# trailing placeholder paragraph that is left over at the end.
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$null = $p7.Range.InsertParagraphAfter()
$p7b = $d.Paragraphs(7)
$null = $p7b.Range.InsertParagraphAfter()

# Clean up paragraph 8 (currently a placeholder with a stray empty run) so it
# collapses to a bare empty paragraph, matching the target "<w:p/>".
$p8 = $d.Paragraphs(8)
$r8 = $p8.Range
$r8.InsertXML((New-FlatOpcXml('<w:p/>')))

# Fill paragraph 7 (no longer the last paragraph) with the new two-run text.
$p7c = $d.Paragraphs(7)
$r7 = $p7c.Range
$xml7 = New-FlatOpcXml('<w:p><w:r><w:t>When looked by individual countr</w:t></w:r><w:r><w:t xml:space="preserve">ies the US has successful campaigns across most days, but the least successful in the shortest and the longest days. Canada has the most success between approximately 225 and 250 days, China between 275 and 300, and Australia around 50 days. All countries had success with shorter and longer campaigns. </w:t></w:r></w:p>')
$r7.InsertXML($xml7)

# Paragraph 9 is now the last paragraph in the body; add one more placeholder
# after it so it is not last, fill it with the new text, then delete the
# spare trailing placeholder paragraph.
$p9 = $d.Paragraphs(9)
$null = $p9.Range.InsertParagraphAfter()

$p9b = $d.Paragraphs(9)
$r9 = $p9b.Range
$xml9 = New-FlatOpcXml('<w:p><w:r><w:t>The length of the campaign is only one factor, other factors should be explored such as the type of play, the time of year in the specific country, the goal, etc.</w:t></w:r></w:p>')
$r9.InsertXML($xml9)

$p10 = $d.Paragraphs(10)
$p10.Range.Delete()
